$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Copy() | Out-Null
$ws.Range("AD1:AF1").PasteSpecial(-4122) | Out-Null

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

for ($r = 2; $r -le 51; $r++) {
    $ws.Cells.Item($r, 30).Value = 58
    $ws.Cells.Item($r, 31).Value = 104
    $ws.Cells.Item($r, 32).Value = 0
}
